$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells H1:J1 (copy style from A1 so border/bold/alignment match) ---
$ws.Range("H1").Value = 'Coverage (raw)'
$ws.Range("I1").Value = 'MatchScore'
$ws.Range("J1").Value = 'Matched'
$ws.Range("A1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)

# --- Row 2 ---
$ws.Range("A2").Value = 'IHA4001YC_13'
$ws.Range("C2").Value = '''1'
$ws.Range("D2").Value = 'ค่าห้องผู้ป่วยปกติ ค่าอาหาร และค่าบริการในโรงพยาบาล'
$ws.Range("E2").Value = 'nan'
$ws.Range("F2").Value = '10,000.00 Per Day / 90 Day Per Disability'
$ws.Range("G2").Value = '4,964,668.30 Per Year'
$ws.Range("H2").Value = 'ค่าห้องผู้ป่วยปกติ ค่าอาหาร และค่าบริการในโรง พยาบาล'
$ws.Range("I2").Value = '''87'
$ws.Range("J2").Value = '''TRUE'

# --- Row 3 ---
$ws.Range("C3").Value = '''2'
$ws.Range("D3").Value = 'ค่าห้องผู้ป่วย ICU ค่าอาหาร และค่าบริการในโรงพยาบาล'
$ws.Range("E3").Value = 'nan'
$ws.Range("F3").Value = '20,000.00 Per Day / 15 Day Per Disability'
$ws.Range("G3").Value = 'nan'
$ws.Range("H3").Value = 'ค่าห้องผู้ป่วย ICU ค่าอาหาร และค่าบริการในโรง พยาบาล'
$ws.Range("I3").Value = '''87'
$ws.Range("J3").Value = '''TRUE'

# --- Row 4 ---
$ws.Range("C4").Value = '''3'
$ws.Range("D4").Value = 'ค่าบริการทางการแพทย์เพื่อการตรวจวินิจฉัย'
$ws.Range("E4").Value = 'nan'
$ws.Range("F4").Value = '150,000.00 Per Disability'
$ws.Range("G4").Value = 'nan'
$ws.Range("H4").Value = 'ค่าบริการทางการแพทย์เพื่อการตรวจวินิจฉัย'
$ws.Range("I4").Value = '''100'
$ws.Range("J4").Value = '''TRUE'

# --- Row 5 ---
$ws.Range("C5").Value = '''4'
$ws.Range("D5").Value = 'ค่าแพทย์ที่ปรึกษาพิเศษ'
$ws.Range("E5").Value = 'nan'
$ws.Range("F5").Value = '10,000.00 Per Disability'
$ws.Range("G5").Value = 'nan'
$ws.Range("H5").Value = 'ค่าแพทย์ที่ปรึกษา'
$ws.Range("I5").Value = '''87'
$ws.Range("J5").Value = '''TRUE'

# --- Row 6 ---
$ws.Range("C6").Value = '''5'
$ws.Range("D6").Value = 'ค่าบริการทางการแพทย์เพื่อการบำบัดรักษา ค่าบริการโลหิตและส่วนประกอบของโลหิต และค่าบริการทางการพยาบาล'
$ws.Range("E6").Value = 'nan'
$ws.Range("F6").Value = '150,000.00 Per Disability'
$ws.Range("G6").Value = 'nan'
$ws.Range("H6").Value = 'ค่าบริการทางการแพทย์เพื่อการบำบัดรักษา ค่าบริการ โลหิตและส่วนประกอบของโลหิต และค่าบริการทางการ พยาบาล'
$ws.Range("I6").Value = '''74'
$ws.Range("J6").Value = '''TRUE'

# --- Row 7 ---
$ws.Range("C7").Value = '''6'
$ws.Range("D7").Value = 'ค่ายา ค่าสารอาหารทางหลอดเลือด และค่าเวชภัณฑ์'
$ws.Range("E7").Value = 'nan'
$ws.Range("F7").Value = '150,000.00 Per Disability'
$ws.Range("G7").Value = 'nan'
$ws.Range("H7").Value = 'ค่ายา ค่าสารอาหารทางหลอดเลือด และค่าวเชกกันท์'
$ws.Range("I7").Value = '''90'
$ws.Range("J7").Value = '''TRUE'

# --- Row 8 ---
$ws.Range("C8").Value = '''7'
$ws.Range("D8").Value = 'ค่ายาและค่าเวชภัณฑ์สิ้นเปลือง (เวชภัณฑ์ 1) สำหรับกลับบ้าน'
$ws.Range("E8").Value = 'nan'
$ws.Range("F8").Value = '150,000.00 Per Disability'
$ws.Range("G8").Value = 'nan'
$ws.Range("H8").Value = 'ค่ายา ค่าเวชภัณฑ์สิ้นเปลือง (เวชภัณฑ์ 1 ) สำหรับกลับ บ้าน'
$ws.Range("I8").Value = '''86'
$ws.Range("J8").Value = '''TRUE'

# --- Row 9 ---
$ws.Range("C9").Value = '''8'
$ws.Range("D9").Value = 'ค่าแพทย์ตรวจรักษา'
$ws.Range("E9").Value = 'nan'
$ws.Range("F9").Value = '1,500.00 Per Day / 1 Visit Per Day / 90 Day Per Disability'
$ws.Range("G9").Value = 'nan'
$ws.Range("H9").Value = 'ค่าแพทย์ตรวจรักษา'
$ws.Range("I9").Value = '''100'
$ws.Range("J9").Value = '''TRUE'

# --- Row 10 ---
$ws.Range("C10").Value = '''9'
$ws.Range("D10").Value = 'ค่าห้องผ่าตัด และค่าห้องทำหัตถการ'
$ws.Range("E10").Value = 'nan'
$ws.Range("F10").Value = '200,000.00 Per Disability'
$ws.Range("G10").Value = 'nan'
$ws.Range("H10").Value = 'ค่าห้องผ้าตัด และค่าห้องทำหัตถการ'
$ws.Range("I10").Value = '''97'
$ws.Range("J10").Value = '''TRUE'

# --- Row 11 ---
$ws.Range("C11").Value = '''10'
$ws.Range("D11").Value = 'ค่ายา ค่าสารอาหารทางหลอดเลือด ค่าเวชภัณฑ์ และค่าอุปกรณ์การผ่าตัดและหัตถการ'
$ws.Range("E11").Value = 'nan'
$ws.Range("F11").Value = '200,000.00 Per Disability'
$ws.Range("G11").Value = 'nan'
$ws.Range("H11").Value = 'ค่ายา ค่าสาธารณางทางหลอดเลือด ค่าเวชภัณฑ์ และค่า อุปกรณ์การฆ่าตัดและหัตถการ'
$ws.Range("I11").Value = '''85'
$ws.Range("J11").Value = '''TRUE'

# --- Row 12 ---
$ws.Range("C12").Value = '''11'
$ws.Range("D12").Value = 'ค่าแพทย์ผ่าตัดและหัตถการ'
$ws.Range("E12").Value = 'nan'
$ws.Range("F12").Value = '200,000.00 Per Disability'
$ws.Range("G12").Value = 'nan'
$ws.Range("H12").Value = 'ค่าแพทย์ผ่าตัดและหัตถการ'
$ws.Range("I12").Value = '''100'
$ws.Range("J12").Value = '''TRUE'

# --- Row 13 ---
$ws.Range("C13").Value = '''12'
$ws.Range("D13").Value = 'ค่าแพทย์ที่ปรึกษา กรณีผ่าตัด (เปลี่ยนอวัยวะ)'
$ws.Range("E13").Value = 'nan'
$ws.Range("F13").Value = '10,000.00 Per Disability'
$ws.Range("G13").Value = 'nan'
$ws.Range("H13").Value = 'ค่าแพทย์ที่ปรึกษา กรณีผ่าตัด'
$ws.Range("I13").Value = '''78'
$ws.Range("J13").Value = '''TRUE'

# --- Row 14 ---
$ws.Range("C14").Value = '''13'
$ws.Range("D14").Value = 'ค่าวิสัญญีแพทย์'
$ws.Range("E14").Value = 'nan'
$ws.Range("F14").Value = '200,000.00 Per Disability'
$ws.Range("G14").Value = 'nan'
$ws.Range("H14").Value = 'ค่าวิสัญญูแพทย์'
$ws.Range("I14").Value = '''93'
$ws.Range("J14").Value = '''TRUE'

# --- Row 15 ---
$ws.Range("C15").Value = '''14'
$ws.Range("D15").Value = 'ค่าห้องผ่าตัด และค่าห้องทำหัตถการ (เปลี่ยนอวัยวะ)'
$ws.Range("E15").Value = 'nan'
$ws.Range("F15").Value = '200,000.00 Per Disability'
$ws.Range("G15").Value = 'nan'
$ws.Range("H15").Value = 'ค่าห้องผ่าตัด และค่าห้องทำหัตถการ (เปลี่ยนอวัยวะ)'
$ws.Range("I15").Value = '''100'
$ws.Range("J15").Value = '''TRUE'

# --- Row 16 ---
$ws.Range("C16").Value = '''15'
$ws.Range("D16").Value = 'ค่ายา ค่าสารอาหารทางหลอดเลือด ค่าเวชภัณฑ์ และค่าอุปกรณ์การผ่าาตัดและหัตถการ (เปลี่ยนอวัยวะ)'
$ws.Range("E16").Value = 'nan'
$ws.Range("F16").Value = '200,000.00 Per Disability'
$ws.Range("G16").Value = 'nan'
$ws.Range("H16").Value = 'ค่ายา ค่าสาธารณางทางหลอดเลือด ค่าเวชภัณฑ์ และค่า อุปกรณ์การฆ่าตัดและหัตถการ (เปลี่ยนอวัยวะ)'
$ws.Range("I16").Value = '''87'
$ws.Range("J16").Value = '''TRUE'

# --- Row 17 ---
$ws.Range("C17").Value = '''16'
$ws.Range("D17").Value = 'ค่าแพทย์ผ่าตัดและหัตถการ (เปลี่ยนอวัยวะ)'
$ws.Range("E17").Value = 'nan'
$ws.Range("F17").Value = '200,000.00 Per Disability'
$ws.Range("G17").Value = 'nan'
$ws.Range("H17").Value = 'ค่าแพทย์ผ่าตัดและหัตถการ (เปลี่ยนอวัยวะ)'
$ws.Range("I17").Value = '''100'
$ws.Range("J17").Value = '''TRUE'

# --- Row 18 ---
$ws.Range("C18").Value = '''17'
$ws.Range("D18").Value = 'ค่าแพทย์ที่ปรึกษา กรณีผ่าตัด (เปลี่ยนอวัยวะ)'
$ws.Range("E18").Value = 'nan'
$ws.Range("F18").Value = '10,000.00 Per Disability'
$ws.Range("G18").Value = 'nan'
$ws.Range("H18").Value = 'ค่าแพทย์ที่ปรึกษา กรณีผ่าตัด (เปลี่ยนอวัยวะ)'
$ws.Range("I18").Value = '''100'
$ws.Range("J18").Value = '''TRUE'

# --- Row 19 ---
$ws.Range("C19").Value = '''18'
$ws.Range("D19").Value = 'ค่าวิสัญญีแพทย์ (เปลี่ยนอวัยวะ)'
$ws.Range("E19").Value = 'nan'
$ws.Range("F19").Value = '200,000.00 Per Disability'
$ws.Range("G19").Value = 'nan'
$ws.Range("H19").Value = 'ค่าวิสัญญูแพทย์ (เปลี่ยนอวัยวะ)'
$ws.Range("I19").Value = '''97'
$ws.Range("J19").Value = '''TRUE'

# --- Row 20 ---
$ws.Range("C20").Value = '''19'
$ws.Range("D20").Value = 'ค่าผ่าตัดใหญ่ (Day Surgery)'
$ws.Range("E20").Value = 'nan'
$ws.Range("F20").Value = '200,000.00 Per Disability'
$ws.Range("G20").Value = 'nan'
$ws.Range("H20").Value = 'ค่าผ่าตัดใหญ่ (Day Surgery)'
$ws.Range("I20").Value = '''100'
$ws.Range("J20").Value = '''TRUE'

# --- Row 21 ---
$ws.Range("C21").Value = '''20'
$ws.Range("D21").Value = 'ค่าบริการทางการแพทย์เพื่อการตรวจวินิฉัยก่อนและหลังจากการรักษาเป็นผู้ป่วยใน'
$ws.Range("E21").Value = 'nan'
$ws.Range("F21").Value = '150,000.00 Per Disability'
$ws.Range("G21").Value = 'nan'
$ws.Range("H21").Value = 'ค่าบริการทางการแพทย์เพื่อการตรวจวินิจฉัยก่อนและ หลังจากการรักษาเป็นผู้ป่วยไข้'
$ws.Range("I21").Value = '''95'
$ws.Range("J21").Value = '''TRUE'

# --- Row 22 ---
$ws.Range("C22").Value = '''21'
$ws.Range("D22").Value = 'ค่ารักษาพยาบาลผู้ป่วยนอกหลังจากการรักษาเป็นผู้ป่วยใน (ไม่รวมค่าบริการทางการแพทย์เพื่อการตรวจวินิจฉัย)'
$ws.Range("E22").Value = 'nan'
$ws.Range("F22").Value = '150,000.00 Per Disability'
$ws.Range("G22").Value = 'nan'
$ws.Range("H22").Value = 'ค่ารักษาพยาบาลผู้ป่วยนอกหลังจากการรักษาเป็นผู้ป่วย ใน (ไม่วรวมค่าบริการทางการแพทย์เพื่อการตรวจวินิจฉัย)'
$ws.Range("I22").Value = '''99'
$ws.Range("J22").Value = '''TRUE'

# --- Row 23 ---
$ws.Range("C23").Value = '''22'
$ws.Range("D23").Value = 'อุบัติเหตุฉุกเฉิน (ER Accident)'
$ws.Range("E23").Value = 'nan'
$ws.Range("F23").Value = '40,000.00 Per Disability'
$ws.Range("G23").Value = 'nan'
$ws.Range("H23").Value = 'อุบัติเหตุฉุกเฉิน (ER Accident)'
$ws.Range("I23").Value = '''100'
$ws.Range("J23").Value = '''TRUE'

# --- Row 24 ---
$ws.Range("C24").Value = '''23'
$ws.Range("D24").Value = 'ค่าเวชศาสตร์ฟื้นฟู หลังการเข้าพักรักษาตัวเป็นผู้ป่วยใน'
$ws.Range("E24").Value = 'nan'
$ws.Range("F24").Value = '150,000.00 Per Disability'
$ws.Range("G24").Value = 'nan'
$ws.Range("H24").Value = 'ค่าวชศสตร์ฟื้นผู้ หลังการเข้ารับการต่อเป็นผู้ป่วยไข้'
$ws.Range("I24").Value = '''81'
$ws.Range("J24").Value = '''TRUE'

# --- Row 25 ---
$ws.Range("C25").Value = '''24'
$ws.Range("D25").Value = 'ค่าบริการทางการแพทย์เพื่อล้างไตผ่านทางเส้นเลือด'
$ws.Range("E25").Value = 'nan'
$ws.Range("F25").Value = '10,000.00 Per Year'
$ws.Range("G25").Value = 'nan'
$ws.Range("H25").Value = 'ค่าบริการทางการแพทย์เพื่อสำอางใดผ่านทางเด้นเลือด'
$ws.Range("I25").Value = '''88'
$ws.Range("J25").Value = '''TRUE'

# --- Row 26 ---
$ws.Range("C26").Value = '''25'
$ws.Range("D26").Value = 'ค่าบริการทางการแพทย์ โดยรังสีรักษา เวชศาสตร์นิวเคลียร์'
$ws.Range("E26").Value = 'nan'
$ws.Range("F26").Value = '10,000.00 Per Year'
$ws.Range("G26").Value = 'nan'
$ws.Range("H26").Value = 'ค่าบริการทางการแพทย์ โดยรังสีรักษา เวชศาสตร์ มิโนคลีทร์'
$ws.Range("I26").Value = '''81'
$ws.Range("J26").Value = '''TRUE'

# --- Row 27 ---
$ws.Range("C27").Value = '''26'
$ws.Range("D27").Value = 'ค่าบริการทางการแพทย์โดยเคมีบำบัด'
$ws.Range("E27").Value = 'nan'
$ws.Range("F27").Value = '10,000.00 Per Year'
$ws.Range("G27").Value = 'nan'
$ws.Range("H27").Value = 'ค่าบริการทางการแพทย์ โดยเคมีบำบัด'
$ws.Range("I27").Value = '''98'
$ws.Range("J27").Value = '''TRUE'

# --- Row 28 ---
$ws.Range("C28").Value = '''27'
$ws.Range("D28").Value = 'ค่าบริการรถพยาบาล'
$ws.Range("E28").Value = 'nan'
$ws.Range("F28").Value = '4,000.00 Per Disability'
$ws.Range("G28").Value = 'nan'
$ws.Range("H28").Value = 'ค่าบริการรถพยาบาล'
$ws.Range("I28").Value = '''100'
$ws.Range("J28").Value = '''TRUE'

# --- Row 29 ---
$ws.Range("C29").Value = '''28'
$ws.Range("D29").Value = 'ค่ารักษาพยาบาล โดยการผ่าตัดเล็ก (Minor Surgery)'
$ws.Range("E29").Value = 'nan'
$ws.Range("F29").Value = '200,000.00 Per Disability'
$ws.Range("G29").Value = 'nan'
$ws.Range("H29").Value = 'ค่ารักษาพยาบาล โดยการผ่าตัดเล็ก (Minor Surgery)'
$ws.Range("I29").Value = '''100'
$ws.Range("J29").Value = '''TRUE'

# --- Remove old trailing rows 30 and 31 (content was folded into earlier rows) ---
$ws.Rows("30:31").Delete()